$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 7642.25
$ws.Range("I64").Value = 7474.75
$ws.Range("J64").Value = 7809.75
$ws.Range("K64").Value = 7474.75
$ws.Range("L64").Value = 7809.75
$ws.Range("M64").Value = -7226.75
$ws.Range("N64").Value = -8305.75

$ws.Range("H67").Value = 7642.25
$ws.Range("I67").Value = 7474.75
$ws.Range("J67").Value = 7809.75
$ws.Range("K67").Value = 7474.75
$ws.Range("L67").Value = 7809.75
$ws.Range("M67").Value = -6616.75
$ws.Range("N67").Value = -9525.75

$ws.Range("H113").Value = 5264.25
$ws.Range("I113").Value = 5936.5557
$ws.Range("J113").Value = 4399.857
$ws.Range("K113").Value = 5936.5557
$ws.Range("L113").Value = 4399.857
$ws.Range("M113").Value = -2682.5557
$ws.Range("N113").Value = -10907.857

$ws.Range("H137").Value = 5871.2593
$ws.Range("I137").Value = 1244.0588
$ws.Range("J137").Value = 13737.5
$ws.Range("K137").Value = 3732.1764
$ws.Range("L137").Value = 41212.5
$ws.Range("M137").Value = -1182.1764
$ws.Range("N137").Value = -46312.5

$ws.Range("H138").Value = 11054.245
$ws.Range("I138").Value = 3044.4
$ws.Range("J138").Value = 13108.052
$ws.Range("K138").Value = 9133.200000000001
$ws.Range("L138").Value = 39324.156
$ws.Range("M138").Value = -3993.200000000001
$ws.Range("N138").Value = -49604.156

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14709223
$ws.Range("I32").Value = 16394967
$ws.Range("K32").Value = 16394967
$ws.Range("M32").Value = -16394680

$ws.Range("H45").Value = 1249.75
$ws.Range("J45").Value = 1500
$ws.Range("L45").Value = 1500
$ws.Range("N45").Value = -2254

$ws.Range("H80").Value = 100000
$ws.Range("J80").Value = 100000
$ws.Range("L80").Value = 100000
$ws.Range("N80").Value = -101996

$ws.Range("H83").Value = 100000
$ws.Range("J83").Value = 100000
$ws.Range("L83").Value = 300000
$ws.Range("N83").Value = -309984

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()

$ws.Range("H105").Value = 9516.200000000001
$ws.Range("I105").Value = 11812.615
$ws.Range("J105").Value = 5251.4287
$ws.Range("K105").Value = 11812.615
$ws.Range("L105").Value = 5251.4287
$ws.Range("M105").Value = -10065.615
$ws.Range("N105").Value = -8745.4287

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 45000
$ws.Range("I4").Value = 45000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 45000
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -44888
$ws.Range("N4").ClearContents()

$ws.Range("H16").Value = 956.46155
$ws.Range("I16").Value = 975.5
$ws.Range("J16").Value = 893
$ws.Range("K16").Value = 975.5
$ws.Range("L16").Value = 893
$ws.Range("M16").Value = -688.5
$ws.Range("N16").Value = -1467

$ws.Range("H31").Value = 21188418
$ws.Range("I31").Value = 1840.5217
$ws.Range("J31").Value = 34724290
$ws.Range("K31").Value = 1840.5217
$ws.Range("L31").Value = 34724290
$ws.Range("M31").Value = -1545.5217
$ws.Range("N31").Value = -34724880

$ws.Range("H34").Value = 21188418
$ws.Range("I34").Value = 1840.5217
$ws.Range("J34").Value = 34724290
$ws.Range("K34").Value = 1840.5217
$ws.Range("L34").Value = 34724290
$ws.Range("M34").Value = -1638.5217
$ws.Range("N34").Value = -34724694

$ws.Range("H113").Value = 956.46155
$ws.Range("I113").Value = 975.5
$ws.Range("J113").Value = 893
$ws.Range("K113").Value = 975.5
$ws.Range("L113").Value = 893
$ws.Range("M113").Value = 1194.5
$ws.Range("N113").Value = -5233

$ws.Range("H122").Value = 2194455
$ws.Range("I122").Value = 1106.7858
$ws.Range("K122").Value = 3320.3574
$ws.Range("M122").Value = -870.3574000000003

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 26977962
$ws.Range("I4").Value = 46820780
$ws.Range("J4").Value = 13371459
$ws.Range("K4").Value = 140462340
$ws.Range("L4").Value = 40114377
$ws.Range("M4").Value = -140462228
$ws.Range("N4").Value = -40114601

$ws.Range("H68").Value = 5329.386
$ws.Range("I68").Value = 5140.2856
$ws.Range("J68").Value = 5365.162
$ws.Range("K68").Value = 15420.8568
$ws.Range("L68").Value = 16095.486
$ws.Range("M68").Value = -14609.8568
$ws.Range("N68").Value = -17717.486

$ws.Range("H71").Value = 5329.386
$ws.Range("I71").Value = 5140.2856
$ws.Range("J71").Value = 5365.162
$ws.Range("K71").Value = 46262.5704
$ws.Range("L71").Value = 48286.458
$ws.Range("M71").Value = -42206.5704
$ws.Range("N71").Value = -56398.458

$ws.Range("H132").Value = 4170181.5
$ws.Range("J132").Value = 6065027.5
$ws.Range("L132").Value = 54585247.5
$ws.Range("N132").Value = -54590307.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

$ws.Range("H7").Value = 3522
$ws.Range("I7").Value = 3580
$ws.Range("K7").Value = 3580
$ws.Range("M7").Value = -3468

$ws.Range("H16").Value = 1350.2858
$ws.Range("I16").Value = 1375.6471
$ws.Range("J16").Value = 1242.5
$ws.Range("K16").Value = 1375.6471
$ws.Range("L16").Value = 1242.5
$ws.Range("M16").Value = -1205.6471
$ws.Range("N16").Value = -1582.5

$ws.Range("H40").Value = 7008.5586
$ws.Range("I40").Value = 6719.161
$ws.Range("J40").Value = 9999
$ws.Range("K40").Value = 6719.161
$ws.Range("L40").Value = 9999
$ws.Range("M40").Value = -6583.161
$ws.Range("N40").Value = -10271

$ws.Range("H81").Value = 56618.285
$ws.Range("I81").Value = 52332
$ws.Range("J81").Value = 62333.332
$ws.Range("K81").Value = 52332
$ws.Range("L81").Value = 62333.332
$ws.Range("M81").Value = -51334
$ws.Range("N81").Value = -64329.332

$ws.Range("H84").Value = 56618.285
$ws.Range("I84").Value = 52332
$ws.Range("J84").Value = 62333.332
$ws.Range("K84").Value = 156996
$ws.Range("L84").Value = 186999.996
$ws.Range("M84").Value = -152004
$ws.Range("N84").Value = -196983.996

$ws.Range("H100").Value = 4410.08
$ws.Range("I100").Value = 3290.8948
$ws.Range("K100").Value = 3290.8948
$ws.Range("M100").Value = -2749.8948

$ws.Range("H122").Value = 2910477.2
$ws.Range("I122").Value = 3203.6177
$ws.Range("J122").Value = 13893511
$ws.Range("K122").Value = 9610.8531
$ws.Range("L122").Value = 41680533
$ws.Range("M122").Value = -7160.8531
$ws.Range("N122").Value = -41685433

$ws.Range("H126").Value = 3522
$ws.Range("I126").Value = 3580
$ws.Range("K126").Value = 10740
$ws.Range("M126").Value = -8270

$ws.Range("H132").Value = 125001704
$ws.Range("I132").Value = 1867.9231
$ws.Range("K132").Value = 5603.7693
$ws.Range("M132").Value = -3073.7693
